$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values replacing the previous "Strike#" totals in column G,
# for rows 2-27 (header row 1 is "K" already).
$kValues = @{
    2  = 5
    3  = 7
    4  = 6
    5  = 2
    6  = 5
    7  = 5
    8  = 3
    9  = 4
    10 = 5
    11 = 1
    12 = 3
    13 = 8
    14 = 5
    15 = 6
    16 = 6
    17 = 3
    18 = 1
    19 = 7
    20 = 6
    21 = 6
    22 = 1
    23 = 2
    24 = 8
    25 = 4
    26 = 5
    27 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
